# Update the division-facts table in place.
# Each data row of the table (rows 1, 5, 9, 13, 17 — the rest are blank
# spacer rows) gets its five cell values replaced with new problems,
# while every other aspect of the document (formatting, row/cell
# structure) is left untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# row -> list of new values for columns 1..5
$updates = @{
    1  = @("33÷8=4, 1", "31÷6=5, 1", "94÷7=13, 3", "41÷8=5, 1", "22÷2=11, 0")
    5  = @("94÷6=15, 4", "97÷2=48, 1", "33÷6=5, 3", "65÷9=7, 2", "60÷3=20, 0")
    9  = @("58÷6=9, 4", "86÷3=28, 2", "53÷9=5, 8", "95÷3=31, 2", "67÷2=33, 1")
    13 = @("63÷9=7, 0", "98÷9=10, 8", "81÷9=9, 0", "58÷3=19, 1", "23÷9=2, 5")
    17 = @("93÷4=23, 1", "53÷2=26, 1", "28÷5=5, 3", "74÷6=12, 2", "72÷8=9, 0")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le $values.Count; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$col - 1]
    }
}
